$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 12 to make room for a new
# "Docentes responsaveis:" (Faculty in charge) section, pushing the
# rest of the sheet (Programa resumido: onward) down by 4 rows.
$ws.Rows.Item(12).Resize(4).Insert()

# The insert carries column A's style into rows 13-15 even though
# those rows only use columns B/C in the final layout - drop them.
$ws.Range("A13:A15").Clear()

# Row 12: new section label, column A only (same look as the other
# bold field labels, e.g. row 16 "Programa resumido:").
$ws.Range("A16").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Rows 13-15: the three faculty members, duplicated in columns B
# (normal text) and C (red text), matching the look of row 16's B/C.
$ws.Range("B16").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

Write-Host "Applied Docentes responsaveis section edit."
